$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 6.044933333333333
$ws.Range("H2").Value = 18.1348
$ws.Range("I2").Value = 0.9708761253868625
$ws.Range("J2").Value = 0.9708761253868624
$ws.Range("M2").Value = 0.978762
$ws.Range("N2").Value = 2.936286
$ws.Range("O2").Value = 0.3819465121442868
$ws.Range("P2").Value = 0.3819465121442868
$ws.Range("Q2").Value = 5.9165510392
$ws.Range("R2").Value = 53.24895935279999
$ws.Range("S2").Value = 0.3708227498156714
$ws.Range("T2").Value = 0.3708227498156714

# Row 3
$ws.Range("G3").Value = 6.044933333333333
$ws.Range("H3").Value = 18.1348
$ws.Range("I3").Value = 0.9708761253868625
$ws.Range("J3").Value = 0.9708761253868624
$ws.Range("O3").Value = 0.0406458950147437
$ws.Range("P3").Value = 0.04064589501474371
$ws.Range("Q3").Value = 0.6296261511555554
$ws.Range("R3").Value = 5.6666353604
$ws.Range("S3").Value = 0.03946212906479556
$ws.Range("T3").Value = 0.03946212906479556

# Row 4
$ws.Range("G4").Value = 6.044933333333333
$ws.Range("H4").Value = 18.1348
$ws.Range("I4").Value = 0.9708761253868625
$ws.Range("J4").Value = 0.9708761253868624
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.3139526666666667
$ws.Range("N4").Value = 0.9418580000000001
$ws.Range("O4").Value = 0.1225151017425393
$ws.Range("P4").Value = 0.1225151017425393
$ws.Range("Q4").Value = 1.897822939822222
$ws.Range("R4").Value = 17.0804064584
$ws.Range("S4").Value = 0.1189469872811738
$ws.Range("T4").Value = 0.1189469872811738

# Row 5
$ws.Range("G5").Value = 6.044933333333333
$ws.Range("H5").Value = 18.1348
$ws.Range("I5").Value = 0.9708761253868625
$ws.Range("J5").Value = 0.9708761253868624
$ws.Range("M5").Value = 0.5829876666666666
$ws.Range("N5").Value = 1.748963
$ws.Range("O5").Value = 0.2275017888991087
$ws.Range("P5").Value = 0.2275017888991087
$ws.Range("Q5").Value = 3.524121579155555
$ws.Range("R5").Value = 31.7170942124
$ws.Range("S5").Value = 0.2208760553249466
$ws.Range("T5").Value = 0.2208760553249466

# Row 6
$ws.Range("G6").Value = 6.044933333333333
$ws.Range("H6").Value = 18.1348
$ws.Range("I6").Value = 0.9708761253868625
$ws.Range("J6").Value = 0.9708761253868624
$ws.Range("M6").Value = 0.3019996666666667
$ws.Range("N6").Value = 0.905999
$ws.Range("O6").Value = 0.117850631054404
$ws.Range("P6").Value = 0.117850631054404
$ws.Range("Q6").Value = 1.825567851688889
$ws.Range("R6").Value = 16.4301106652
$ws.Range("S6").Value = 0.1144183640524964
$ws.Range("T6").Value = 0.1144183640524964

# Row 7
$ws.Range("G7").Value = 6.044933333333333
$ws.Range("H7").Value = 18.1348
$ws.Range("I7").Value = 0.9708761253868625
$ws.Range("J7").Value = 0.9708761253868624
$ws.Range("M7").Value = 0.2807033333333333
$ws.Range("N7").Value = 0.8421099999999999
$ws.Range("O7").Value = 0.1095400711449175
$ws.Range("P7").Value = 0.1095400711449176
$ws.Range("Q7").Value = 1.696832936444444
$ws.Range("R7").Value = 15.271496428
$ws.Range("S7").Value = 0.1063498398477788
$ws.Range("T7").Value = 0.1063498398477788

# Row 8
$ws.Range("I8").Value = 0.01821359071319307
$ws.Range("J8").Value = 0.01821359071319307
$ws.Range("M8").Value = 0.978762
$ws.Range("N8").Value = 2.936286
$ws.Range("O8").Value = 0.3819465121442868
$ws.Range("P8").Value = 0.3819465121442868
$ws.Range("Q8").Value = 0.110994220832
$ws.Range("R8").Value = 0.9989479874880001
$ws.Range("S8").Value = 0.006956617446527667
$ws.Range("T8").Value = 0.006956617446527667

# Row 9
$ws.Range("I9").Value = 0.01821359071319307
$ws.Range("J9").Value = 0.01821359071319307
$ws.Range("O9").Value = 0.0406458950147437
$ws.Range("P9").Value = 0.04064589501474371
$ws.Range("S9").Value = 0.0007403076959699564
$ws.Range("T9").Value = 0.0007403076959699566

# Row 10
$ws.Range("I10").Value = 0.01821359071319307
$ws.Range("J10").Value = 0.01821359071319307
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.3139526666666667
$ws.Range("N10").Value = 0.9418580000000001
$ws.Range("O10").Value = 0.1225151017425393
$ws.Range("P10").Value = 0.1225151017425393
$ws.Range("Q10").Value = 0.03560306960711112
$ws.Range("R10").Value = 0.320427626464
$ws.Range("S10").Value = 0.002231439919323818
$ws.Range("T10").Value = 0.002231439919323818

# Row 11
$ws.Range("I11").Value = 0.01821359071319307
$ws.Range("J11").Value = 0.01821359071319307
$ws.Range("M11").Value = 0.5829876666666666
$ws.Range("N11").Value = 1.748963
$ws.Range("O11").Value = 0.2275017888991087
$ws.Range("P11").Value = 0.2275017888991087
$ws.Range("Q11").Value = 0.06611235603377777
$ws.Range("R11").Value = 0.595011204304
$ws.Range("S11").Value = 0.004143624469527616
$ws.Range("T11").Value = 0.004143624469527618

# Row 12
$ws.Range("I12").Value = 0.01821359071319307
$ws.Range("J12").Value = 0.01821359071319307
$ws.Range("M12").Value = 0.3019996666666667
$ws.Range("N12").Value = 0.905999
$ws.Range("O12").Value = 0.117850631054404
$ws.Range("P12").Value = 0.117850631054404
$ws.Range("Q12").Value = 0.03424756753244444
$ws.Range("R12").Value = 0.308228107792
$ws.Range("S12").Value = 0.002146483159316436
$ws.Range("T12").Value = 0.002146483159316436

# Row 13
$ws.Range("I13").Value = 0.01821359071319307
$ws.Range("J13").Value = 0.01821359071319307
$ws.Range("M13").Value = 0.2807033333333333
$ws.Range("N13").Value = 0.8421099999999999
$ws.Range("O13").Value = 0.1095400711449175
$ws.Range("P13").Value = 0.1095400711449176
$ws.Range("Q13").Value = 0.03183250654222222
$ws.Range("R13").Value = 0.28649255888
$ws.Range("S13").Value = 0.001995118022527578
$ws.Range("T13").Value = 0.001995118022527579

# Row 14
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 0.3333333333333333
$ws.Range("G14").Value = 0.06793033333333333
$ws.Range("H14").Value = 0.203791
$ws.Range("I14").Value = 0.01091028389994453
$ws.Range("J14").Value = 0.01091028389994453
$ws.Range("M14").Value = 0.978762
$ws.Range("N14").Value = 2.936286
$ws.Range("O14").Value = 0.3819465121442868
$ws.Range("P14").Value = 0.3819465121442868
$ws.Range("Q14").Value = 0.06648762891399999
$ws.Range("R14").Value = 0.598388660226
$ws.Range("S14").Value = 0.004167144882087782
$ws.Range("T14").Value = 0.004167144882087782

# Row 15
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0.3333333333333333
$ws.Range("G15").Value = 0.06793033333333333
$ws.Range("H15").Value = 0.203791
$ws.Range("I15").Value = 0.01091028389994453
$ws.Range("J15").Value = 0.01091028389994453
$ws.Range("O15").Value = 0.0406458950147437
$ws.Range("P15").Value = 0.04064589501474371
$ws.Range("Q15").Value = 0.007075465015888888
$ws.Range("R15").Value = 0.063679185143
$ws.Range("S15").Value = 0.0004434582539781939
$ws.Range("T15").Value = 0.000443458253978194

# Row 16
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0.3333333333333333
$ws.Range("G16").Value = 0.06793033333333333
$ws.Range("H16").Value = 0.203791
$ws.Range("I16").Value = 0.01091028389994453
$ws.Range("J16").Value = 0.01091028389994453
$ws.Range("M16").Value = 0.3139526666666667
$ws.Range("N16").Value = 0.9418580000000001
$ws.Range("O16").Value = 0.1225151017425393
$ws.Range("P16").Value = 0.1225151017425393
$ws.Range("Q16").Value = 0.02132690929755556
$ws.Range("R16").Value = 0.191942183678
$ws.Range("S16").Value = 0.001336674542041693
$ws.Range("T16").Value = 0.001336674542041693

# Row 17
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 0.3333333333333333
$ws.Range("G17").Value = 0.06793033333333333
$ws.Range("H17").Value = 0.203791
$ws.Range("I17").Value = 0.01091028389994453
$ws.Range("J17").Value = 0.01091028389994453
$ws.Range("M17").Value = 0.5829876666666666
$ws.Range("N17").Value = 1.748963
$ws.Range("O17").Value = 0.2275017888991087
$ws.Range("P17").Value = 0.2275017888991087
$ws.Range("Q17").Value = 0.03960254652588888
$ws.Range("R17").Value = 0.356422918733
$ws.Range("S17").Value = 0.002482109104634525
$ws.Range("T17").Value = 0.002482109104634526

# Row 18
$ws.Range("E18").Value = 1
$ws.Range("F18").Value = 0.3333333333333333
$ws.Range("G18").Value = 0.06793033333333333
$ws.Range("H18").Value = 0.203791
$ws.Range("I18").Value = 0.01091028389994453
$ws.Range("J18").Value = 0.01091028389994453
$ws.Range("M18").Value = 0.3019996666666667
$ws.Range("N18").Value = 0.905999
$ws.Range("O18").Value = 0.117850631054404
$ws.Range("P18").Value = 0.117850631054404
$ws.Range("Q18").Value = 0.02051493802322222
$ws.Range("R18").Value = 0.184634442209
$ws.Range("S18").Value = 0.001285783842591167
$ws.Range("T18").Value = 0.001285783842591167

# Row 19
$ws.Range("E19").Value = 1
$ws.Range("F19").Value = 0.3333333333333333
$ws.Range("G19").Value = 0.06793033333333333
$ws.Range("H19").Value = 0.203791
$ws.Range("I19").Value = 0.01091028389994453
$ws.Range("J19").Value = 0.01091028389994453
$ws.Range("M19").Value = 0.2807033333333333
$ws.Range("N19").Value = 0.8421099999999999
$ws.Range("O19").Value = 0.1095400711449175
$ws.Range("P19").Value = 0.1095400711449176
$ws.Range("Q19").Value = 0.01906827100111111
$ws.Range("R19").Value = 0.17161443901
$ws.Range("S19").Value = 0.001195113274611172
$ws.Range("T19").Value = 0.001195113274611173
